# Fruta / hortaliza, semanal
#
# A new weekly price observation was inserted as row 117 of the "Choclo"
# sheet (Comercializadora del Agro de Limarí). All the rows that used to
# occupy 117-139 shift down by one, to 118-140, and the sheet's used range
# grows from A1:R139 to A1:R140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 117:139 down to 118:140, leaving a blank row 117.
$ws.Rows("117:117").Insert()

# Populate the newly-inserted row 117 with the new observation.
$ws.Range("A117").Value = 2
$ws.Range("B117").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C117").Value = 'Coquimbo'
$ws.Range("D117").Value = 44889
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100112024
$ws.Range("G117").Value = 'Choclo'
$ws.Range("H117").Value = 'Dulce o Americano'
$ws.Range("I117").Value = 'Primera'
$ws.Range("J117").Value = 460
$ws.Range("K117").Value = 26000
$ws.Range("L117").Value = 27000
$ws.Range("M117").Value = 26500
$ws.Range("N117").Value = '$/malla 70 unidades'
$ws.Range("O117").Value = 'Provincia de Limarí'
$ws.Range("P117").Value = 379
$ws.Range("Q117").Value = 70
$ws.Range("R117").Value = 'Hortaliza'
